$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/volume table refresh (GitHub Actions scheduled update).
# Price cells that look like plain numbers ("1.00", "7.79", ...) must stay
# TEXT (the sheet stores prices as formatted strings, not numeric cells), so
# those are pre-formatted as Text before the value is written; this keeps
# Excel from silently re-parsing "0.580" -> 0.58, "1.00" -> 1, etc.

$ws.Range("D2").Value = '50.942.37'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '2.945.53'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '378.92'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.08'
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.07'
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0848'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.79'
$ws.Range("E13").Value = '  +4.52%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.405.39'
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.26'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("E16").Value = '  +68.81%  '
$ws.Range("D17").Value = '2.951.78'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").Value = '50.915.50'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("E20").Value = '  -4.51%  '
$ws.Range("E21").Value = '  -1.92%  '
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.34'
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.15'
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("E25").Value = '  +9.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.10'
$ws.Range("E26").Value = '  -3.29%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.04'
$ws.Range("E28").Value = '  -10.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.55'
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("E30").Value = '  -3.99%  '
$ws.Range("E31").Value = '  -4.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.09'
$ws.Range("E32").Value = '  +2.50%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.44'
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.05'
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.43'
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0429'
$ws.Range("E36").Value = '  -6.11%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.08'
$ws.Range("E38").Value = '  +2.66%  '
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.52'
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.50'
$ws.Range("E42").Value = '  -2.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '117.77'
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.51'
$ws.Range("E44").Value = '  +7.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.24'
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("E47").Value = '  -2.25%  '
$ws.Range("D48").Value = '1.997.87'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.259'
$ws.Range("E49").Value = '  -4.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0314'
$ws.Range("E50").Value = '  -10.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.27'
$ws.Range("E51").Value = '  +3.56%  '
